$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.581.56"
$ws.Range("E2").Value = "  +2.23%  "

$ws.Range("D3").Value = "2.061.18"
$ws.Range("E3").Value = "  +8.95%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.87"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  -4.29%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.21"
$ws.Range("E8").Value = "  +4.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.59"
$ws.Range("E9").Value = "  +6.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.363"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0722"
$ws.Range("E11").Value = "  -4.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0985"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").Value = "2.360.00"
$ws.Range("E14").Value = "  +9.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.823"
$ws.Range("E15").Value = "  +4.04%  "

$ws.Range("D16").Value = "2.063.07"
$ws.Range("E16").Value = "  +9.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.92"
$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").Value = "36.586.82"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.37"
$ws.Range("E19").Value = "  -3.02%  "

$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "238.30"
$ws.Range("E21").Value = "  -3.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.66"
$ws.Range("E22").Value = "  -3.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.90"
$ws.Range("E23").Value = "  -5.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("E25").Value = "  -7.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.69"
$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.88"
$ws.Range("E27").Value = "  +2.45%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.27"
$ws.Range("E28").Value = "  +10.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("E29").Value = "  -9.53%  "

$ws.Range("E30").Value = "  -5.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.38"
$ws.Range("E31").Value = "  +50.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.35"
$ws.Range("E32").Value = "  -1.70%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0922"
$ws.Range("E33").Value = "  +23.46%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0581"
$ws.Range("E34").Value = "  -4.22%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.34"
$ws.Range("E35").Value = "  +20.64%  "

$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  -6.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.886"
$ws.Range("E39").Value = "  +3.48%  "

$ws.Range("E40").Value = "  -10.76%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.01"
$ws.Range("E41").Value = "  -2.11%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("E42").Value = "  -6.04%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.12"
$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("E44").Value = "  +15.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.89"
$ws.Range("E45").Value = "  -6.07%  "

$ws.Range("D46").Value = "1.317.71"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0820"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").Value = "2.254.61"
$ws.Range("E49").Value = "  +9.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -5.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.84"
$ws.Range("E51").Value = "  +14.50%  "
